$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.413.85"
$ws.Range("E2").Value = "  +5.17%  "
$ws.Range("D3").Value = "1.814.24"
$ws.Range("E3").Value = "  +5.02%  "
$ws.Range("D4").Value = "'0.9971"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'317.70"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.5650"
$ws.Range("E7").Value = "  +16.42%  "
$ws.Range("D8").Value = "'0.3835"
$ws.Range("E8").Value = "  +9.93%  "
$ws.Range("D9").Value = "'43.38"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'0.07624"
$ws.Range("E10").Value = "  +5.24%  "
$ws.Range("E11").Value = "  +7.80%  "
$ws.Range("D12").Value = "'21.33"
$ws.Range("E12").Value = "  +6.81%  "
$ws.Range("D13").Value = "'0.9964"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +6.07%  "
$ws.Range("D15").Value = "1.800.33"
$ws.Range("E15").Value = "  +4.31%  "
$ws.Range("D16").Value = "'7.233"
$ws.Range("E16").Value = "  +5.49%  "
$ws.Range("D17").Value = "'92.27"
$ws.Range("E17").Value = "  +6.08%  "
$ws.Range("D18").Value = "'0.00001079"
$ws.Range("E18").Value = "  +4.32%  "
$ws.Range("D19").Value = "'0.06505"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'0.9961"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").Value = "'5.991"
$ws.Range("E22").Value = "  +4.82%  "
$ws.Range("D23").Value = "28.418.91"
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("D24").Value = "'11.29"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "'2.108"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").Value = "'156.78"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").Value = "'2.379"
$ws.Range("E28").Value = "  +14.48%  "
$ws.Range("D29").Value = "2.013.79"
$ws.Range("E29").Value = "  +4.74%  "
$ws.Range("D30").Value = "'123.57"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "'1.145"
$ws.Range("E31").Value = "  +9.40%  "
$ws.Range("D32").Value = "'0.1046"
$ws.Range("E32").Value = "  +12.05%  "
$ws.Range("D33").Value = "'5.764"
$ws.Range("E33").Value = "  +6.73%  "
$ws.Range("D34").Value = "'3.628"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("D36").Value = "'0.2139"
$ws.Range("E36").Value = "  +7.13%  "
$ws.Range("D37").Value = "'8.703"
$ws.Range("E37").Value = "  +15.32%  "
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06078"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").Value = "'5.042"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6409"
$ws.Range("E41").Value = "  +6.95%  "
$ws.Range("D42").Value = "'0.9954"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "'1.153"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").Value = "'1.377"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = "  +5.87%  "
$ws.Range("D46").Value = "'0.5991"
$ws.Range("E46").Value = "  +6.42%  "
$ws.Range("D47").Value = "'3.703"
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("D48").Value = "'122.60"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "'1.936"
$ws.Range("E49").Value = "  +4.70%  "
$ws.Range("D50").Value = "'1.144"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").Value = "'0.06833"
$ws.Range("E51").Value = "  +2.96%  "
